$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("expected")

# Set A3 on the "expected" sheet to the new label "index"
$ws2.Range("A3").Value = "index"

# Clear the now-redundant "Other" labels in B5 and B8 (keep formatting/style intact)
$ws2.Range("B5").ClearContents()
$ws2.Range("B8").ClearContents()

# Make "expected" the active sheet/tab so it becomes tabSelected and the workbook's activeTab
$ws2.Activate()

# Update the selection on the "expected" sheet to A4
$ws2.Range("A4").Select()

$wb.Save()
